# "Factory Method" slide (slide 3) - the description bullet currently ends
# with "... cho lop con (subclass)"; it should end with
# "... cho lop con (subclass)." (a trailing period is added).
$p = $ppt.ActivePresentation

$targetSlideIndex = 3
$targetShapeIndex = 2
$needle = "subclass)"

$s = $p.Slides.Item($targetSlideIndex)
$shape = $s.Shapes.Item($targetShapeIndex)
$tr = $shape.TextFrame.TextRange
$full = $tr.Text

$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $len = $full.Length - $idx
    $target = $tr.Characters($idx + 1, $len)
    if ($target.Text -eq $needle) {
        $target.Text = $needle + "."
    }
}
